# Scheduled runner: refresh market-board derived columns (H:N) on the
# per-job Leve profit sheets. Values recomputed from latest price pulls.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 198.33333
$ws.Range("I9").Value = 178
$ws.Range("J9").Value = 300
$ws.Range("K9").Value = 178
$ws.Range("L9").Value = 300
$ws.Range("M9").Value = -9
$ws.Range("N9").Value = -638

$ws.Range("H92").Value = 728.5714
$ws.Range("I92").Value = 771.25
$ws.Range("K92").Value = 771.25
$ws.Range("M92").Value = 476.75

$ws.Range("H129").Value = 830
$ws.Range("I129").Value = 483.33334
$ws.Range("J129").Value = 1003.3333
$ws.Range("K129").Value = 1450.00002
$ws.Range("L129").Value = 3009.9999
$ws.Range("M129").Value = 3549.99998
$ws.Range("N129").Value = -13009.9999

$ws.Range("H138").Value = 2608.161
$ws.Range("J138").Value = 2791.697
$ws.Range("L138").Value = 8375.091
$ws.Range("N138").Value = -18655.091


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 982138.4399999999
$ws.Range("I2").Value = 1565.5
$ws.Range("J2").Value = 2452998
$ws.Range("K2").Value = 1565.5
$ws.Range("L2").Value = 2452998
$ws.Range("M2").Value = -1452.5
$ws.Range("N2").Value = -2453224

$ws.Range("H5").Value = 63.333332
$ws.Range("I5").Value = 63.333332
$ws.Range("K5").Value = 63.333332
$ws.Range("M5").Value = 48.666668

$ws.Range("H63").Value = 2850.5264
$ws.Range("I63").Value = 1505.4546
$ws.Range("J63").Value = 4700
$ws.Range("K63").Value = 1505.4546
$ws.Range("L63").Value = 4700
$ws.Range("M63").Value = -819.4546
$ws.Range("N63").Value = -6072

$ws.Range("H66").Value = 2850.5264
$ws.Range("I66").Value = 1505.4546
$ws.Range("J66").Value = 4700
$ws.Range("K66").Value = 7527.273
$ws.Range("L66").Value = 23500
$ws.Range("M66").Value = -4095.273
$ws.Range("N66").Value = -30364

$ws.Range("H116").Value = 982138.4399999999
$ws.Range("I116").Value = 1565.5
$ws.Range("J116").Value = 2452998
$ws.Range("K116").Value = 1565.5
$ws.Range("L116").Value = 2452998
$ws.Range("M116").Value = 728.5
$ws.Range("N116").Value = -2457586


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 982138.4399999999
$ws.Range("I3").Value = 1565.5
$ws.Range("J3").Value = 2452998
$ws.Range("K3").Value = 1565.5
$ws.Range("L3").Value = 2452998
$ws.Range("M3").Value = -1451.5
$ws.Range("N3").Value = -2453226

$ws.Range("H4").Value = 63.333332
$ws.Range("I4").Value = 63.333332
$ws.Range("K4").Value = 63.333332
$ws.Range("M4").Value = 51.666668

$ws.Range("H86").Value = 942.08
$ws.Range("I86").Value = 958.2778
$ws.Range("J86").Value = 900.4286
$ws.Range("K86").Value = 958.2778
$ws.Range("L86").Value = 900.4286
$ws.Range("M86").Value = 164.7222
$ws.Range("N86").Value = -3146.4286

$ws.Range("H89").Value = 942.08
$ws.Range("I89").Value = 958.2778
$ws.Range("J89").Value = 900.4286
$ws.Range("K89").Value = 4791.389
$ws.Range("L89").Value = 4502.143
$ws.Range("M89").Value = 824.6109999999999
$ws.Range("N89").Value = -15734.143

$ws.Range("H97").Value = 13000
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 13000
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 13000
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -14982


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 22916.639
$ws.Range("I99").Value = 1850.1786
$ws.Range("J99").Value = 53961.95
$ws.Range("K99").Value = 1850.1786
$ws.Range("L99").Value = 53961.95
$ws.Range("M99").Value = -352.1786
$ws.Range("N99").Value = -56957.95

$ws.Range("H126").Value = 22916.639
$ws.Range("I126").Value = 1850.1786
$ws.Range("J126").Value = 53961.95
$ws.Range("K126").Value = 5550.5358
$ws.Range("L126").Value = 161885.85
$ws.Range("M126").Value = -3080.5358
$ws.Range("N126").Value = -166825.85


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 50133.45
$ws.Range("I12").Value = 85
$ws.Range("J12").Value = 91082.17999999999
$ws.Range("K12").Value = 255
$ws.Range("L12").Value = 273246.54
$ws.Range("M12").Value = -82
$ws.Range("N12").Value = -273592.54

$ws.Range("H37").Value = 32128.072
$ws.Range("J37").Value = 32128.072
$ws.Range("L37").Value = 96384.216
$ws.Range("N37").Value = -96608.216

$ws.Range("H140").Value = 1412.7273
$ws.Range("I140").Value = 1412.7273
$ws.Range("K140").Value = 4238.1819
$ws.Range("M140").Value = 941.8181000000004


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 3001.9443
$ws.Range("I36").Value = 348.2
$ws.Range("J36").Value = 4022.6155
$ws.Range("K36").Value = 348.2
$ws.Range("L36").Value = 4022.6155
$ws.Range("M36").Value = 136.8
$ws.Range("N36").Value = -4992.6155

$ws.Range("H92").Value = 8633.333000000001
$ws.Range("J92").Value = 8633.333000000001
$ws.Range("L92").Value = 8633.333000000001
$ws.Range("N92").Value = -12377.333


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1419.375
$ws.Range("I61").Value = 1410.7142
$ws.Range("J61").Value = 1480
$ws.Range("K61").Value = 1410.7142
$ws.Range("L61").Value = 1480
$ws.Range("M61").Value = -1208.7142
$ws.Range("N61").Value = -1884

$ws.Range("H113").Value = 1419.375
$ws.Range("I113").Value = 1410.7142
$ws.Range("J113").Value = 1480
$ws.Range("K113").Value = 1410.7142
$ws.Range("L113").Value = 1480
$ws.Range("M113").Value = 759.2858000000001
$ws.Range("N113").Value = -5820

$ws.Range("H132").Value = 8711.471
$ws.Range("I132").Value = 9506.929
$ws.Range("J132").Value = 4999.3335
$ws.Range("K132").Value = 28520.787
$ws.Range("L132").Value = 14998.0005
$ws.Range("M132").Value = -25990.787
$ws.Range("N132").Value = -20058.0005


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4909
$ws.Range("I81").Value = 2373.5
$ws.Range("J81").Value = 9980
$ws.Range("K81").Value = 4747
$ws.Range("L81").Value = 19960
$ws.Range("M81").Value = -3686
$ws.Range("N81").Value = -22082

$ws.Range("H84").Value = 4909
$ws.Range("I84").Value = 2373.5
$ws.Range("J84").Value = 9980
$ws.Range("K84").Value = 23735
$ws.Range("L84").Value = 99800
$ws.Range("M84").Value = -18431
$ws.Range("N84").Value = -110408

$ws.Range("H107").Value = 498.8
$ws.Range("I107").Value = 304.3
$ws.Range("J107").Value = 693.3
$ws.Range("K107").Value = 912.9000000000001
$ws.Range("L107").Value = 2079.9
$ws.Range("M107").Value = 1007.1
$ws.Range("N107").Value = -5919.9

$ws.Range("H136").Value = 7038.45
$ws.Range("I136").Value = 9424.393
$ws.Range("J136").Value = 1471.25
$ws.Range("K136").Value = 28273.179
$ws.Range("L136").Value = 4413.75
$ws.Range("M136").Value = -25723.179
$ws.Range("N136").Value = -9513.75
